$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the previously empty row 31 with a new time-tracking entry
$ws.Range("A31").Value = 44055
$ws.Range("B31").Value = 2.5
$ws.Range("C31").Value = "Logon suunnittelu, kuvan lataamisen ja optimointi"

# The added text wraps onto multiple lines, so the row is taller
$ws.Rows.Item(31).RowHeight = 30

# Move the active selection to the newly filled description cell
$ws.Range("C31").Select()
